{"js": "// Update the date line and the 25 division problems in the 5x5 practice\n// table. Cell targets are addressed by (row, col) position rather than by\n// searching for the old text, because several of the new values collide\n// with other *old* values elsewhere in the table (e.g. \"61\u00f77=\" is both a\n// pre-edit cell value and a post-edit replacement for a different cell),\n// which would make a naive global text-replace double-apply.\n\nconst body = context.document.body;\n\n// --- Title / date paragraph (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2025-02-02 Sunday\", \"Replace\");\n\n// --- Division problems table (5 data rows x 5 columns) ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, in row-major order matching the existing 5x5 data rows.\nconst newValues = [\n  [\"58\u00f74=\", \"21\u00f76=\", \"93\u00f78=\", \"20\u00f76=\", \"83\u00f78=\"],\n  [\"17\u00f74=\", \"59\u00f74=\", \"27\u00f74=\", \"90\u00f76=\", \"61\u00f77=\"],\n  [\"34\u00f72=\", \"79\u00f78=\", \"11\u00f72=\", \"52\u00f76=\", \"58\u00f77=\"],\n  [\"96\u00f73=\", \"78\u00f73=\", \"41\u00f75=\", \"28\u00f77=\", \"51\u00f75=\"],\n  [\"15\u00f74=\", \"73\u00f73=\", \"67\u00f73=\", \"19\u00f79=\", \"52\u00f73=\"],\n];\n\n// The table has 20 rows total (5 data rows, each followed by 3 blank\n// rows). Data rows therefore live at table-row indices 0, 4, 8, 12, 16.\nconst dataRowIndices = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < dataRowIndices.length; r++) {\n  const tableRow = dataRowIndices[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(tableRow, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the 5x5 practice\n# table. Cells are addressed by their (row, column) position rather than by\n# searching for the old text, because several of the new values collide\n# with other *old* values elsewhere in the table (e.g. \"61\u00f77=\" is both a\n# pre-edit cell value and a post-edit replacement for a different cell),\n# which would make a naive global find/replace double-apply.\n\n$d = $word.ActiveDocument\n\n# --- Title / date paragraph (first paragraph in the document) ---\n$d.Paragraphs.Item(1).Range.Text = \"2025-02-02 Sunday\"\n\n# --- Division problems table (5 data rows x 5 columns) ---\n$t = $d.Tables.Item(1)\n\n# New values, in row-major order matching the table's existing 5 data rows.\n$newValues = @(\n    @(\"58\u00f74=\", \"21\u00f76=\", \"93\u00f78=\", \"20\u00f76=\", \"83\u00f78=\"),\n    @(\"17\u00f74=\", \"59\u00f74=\", \"27\u00f74=\", \"90\u00f76=\", \"61\u00f77=\"),\n    @(\"34\u00f72=\", \"79\u00f78=\", \"11\u00f72=\", \"52\u00f76=\", \"58\u00f77=\"),\n    @(\"96\u00f73=\", \"78\u00f73=\", \"41\u00f75=\", \"28\u00f77=\", \"51\u00f75=\"),\n    @(\"15\u00f74=\", \"73\u00f73=\", \"67\u00f73=\", \"19\u00f79=\", \"52\u00f73=\")\n)\n\n# The table has 20 rows total (5 data rows, each followed by 3 blank rows).\n# Data rows therefore live at 1-based Word row indices 1, 5, 9, 13, 17.\n$dataRowIndices = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $dataRowIndices.Length; $r++) {\n    $tableRow = $dataRowIndices[$r]\n    for ($c = 0; $c -lt $newValues[$r].Length; $c++) {\n        $t.Cell($tableRow, $c + 1).Range.Text = $newValues[$r][$c]\n    }\n}\n"}
